$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Pipeline(steps=[(''scaler'', StandardScaler()), (''selector'', None),
                (''model'',
                 <class ''skorch.classifier.NeuralNetClassifier''>[initialized](
  module_=TorchMLPClassifier(
    (fc1): Linear(in_features=40, out_features=4, bias=True)
    (relu): ReLU()
    (fc2): Linear(in_features=4, out_features=1, bias=True)
    (softmax): Softmax(dim=1)
    (sigmoid): Sigmoid()
  ),
))])'
$ws.Range("C2").Value = '{''selector'': None, ''scaler'': StandardScaler(), ''model__optimizer__lr'': 0.5, ''model__optimizer'': <class ''torch.optim.adam.Adam''>, ''model__module__hidden_size'': 4, ''model__max_epochs'': 200}'
$ws.Range("D2").Value = 0.5838280809243253
$ws.Range("F2").Value = 0.6793974437005478
$ws.Range("G2").Value = 0.5440579710144928
$ws.Range("I2").Value = '[0 1 1 1 0 0 1 0 1 1 1 1 0 1 1 0 1 0 0 0 0 0 0 0]'

$ws.Range("B3").Value = 'Pipeline(steps=[(''scaler'', StandardScaler()), (''selector'', None),
                (''model'',
                 <class ''skorch.classifier.NeuralNetClassifier''>[initialized](
  module_=TorchMLPClassifier(
    (fc1): Linear(in_features=40, out_features=3, bias=True)
    (relu): ReLU()
    (fc2): Linear(in_features=3, out_features=1, bias=True)
    (softmax): Softmax(dim=1)
    (sigmoid): Sigmoid()
  ),
))])'
$ws.Range("C3").Value = '{''selector'': None, ''scaler'': StandardScaler(), ''model__optimizer__lr'': 0.6, ''model__optimizer'': <class ''torch.optim.adam.Adam''>, ''model__module__hidden_size'': 3, ''model__max_epochs'': 500}'
$ws.Range("D3").Value = 0.6057850245654794
$ws.Range("F3").Value = 0.5295912271718723
$ws.Range("G3").Value = 0.4298245614035088
$ws.Range("I3").Value = '[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]'

$ws.Range("B4").Value = 'Pipeline(steps=[(''scaler'', StandardScaler()), (''selector'', None),
                (''model'',
                 <class ''skorch.classifier.NeuralNetClassifier''>[initialized](
  module_=TorchMLPClassifier(
    (fc1): Linear(in_features=40, out_features=1, bias=True)
    (relu): ReLU()
    (fc2): Linear(in_features=1, out_features=1, bias=True)
    (softmax): Softmax(dim=1)
    (sigmoid): Sigmoid()
  ),
))])'
$ws.Range("C4").Value = '{''selector'': None, ''scaler'': StandardScaler(), ''model__optimizer__lr'': 0.5, ''model__optimizer'': <class ''torch.optim.adam.Adam''>, ''model__module__hidden_size'': 1, ''model__max_epochs'': 400}'
$ws.Range("D4").Value = 0.5588421136453174
$ws.Range("F4").Value = 0.6692209450830141
$ws.Range("G4").Value = 0.5862470862470862
$ws.Range("I4").Value = '[1 0 0 0 1 1 0 1 1 1 1 0 0 1 1 1 0 1 0 0 0 0 1 0]'

$ws.Range("B5").Value = 'Pipeline(steps=[(''scaler'', StandardScaler()), (''selector'', None),
                (''model'',
                 <class ''skorch.classifier.NeuralNetClassifier''>[initialized](
  module_=TorchMLPClassifier(
    (fc1): Linear(in_features=40, out_features=5, bias=True)
    (relu): ReLU()
    (fc2): Linear(in_features=5, out_features=1, bias=True)
    (softmax): Softmax(dim=1)
    (sigmoid): Sigmoid()
  ),
))])'
$ws.Range("C5").Value = '{''selector'': None, ''scaler'': StandardScaler(), ''model__optimizer__lr'': 0.7, ''model__optimizer'': <class ''torch.optim.adam.Adam''>, ''model__module__hidden_size'': 5, ''model__max_epochs'': 400}'
$ws.Range("D5").Value = 0.5708510950021247
$ws.Range("F5").Value = 0.685137924119626
$ws.Range("G5").Value = 0.4895104895104894
$ws.Range("I5").Value = '[0 0 1 0 0 1 1 0 0 0 0 0 1 1 0 0 0 1 1 0 0 1 0 0]'

$ws.Range("B6").Value = 'Pipeline(steps=[(''scaler'', StandardScaler()), (''selector'', None),
                (''model'',
                 <class ''skorch.classifier.NeuralNetClassifier''>[initialized](
  module_=TorchMLPClassifier(
    (fc1): Linear(in_features=40, out_features=3, bias=True)
    (relu): ReLU()
    (fc2): Linear(in_features=3, out_features=1, bias=True)
    (softmax): Softmax(dim=1)
    (sigmoid): Sigmoid()
  ),
))])'
$ws.Range("C6").Value = '{''selector'': None, ''scaler'': StandardScaler(), ''model__optimizer__lr'': 0.7, ''model__optimizer'': <class ''torch.optim.adam.Adam''>, ''model__module__hidden_size'': 3, ''model__max_epochs'': 300}'
$ws.Range("D6").Value = 0.6064671482628139
$ws.Range("F6").Value = 0.6905615292712067
$ws.Range("G6").Value = 0.6269565217391304
$ws.Range("I6").Value = '[0 0 1 1 1 0 0 0 1 0 1 1 0 1 1 1 0 0 0 1 0 0 0 1]'
